$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
# B1: "Correo electrónico" -> "Email"
$ws.Range("B1").Value = "Email"

# D1/E1/F1 text stays the same (Latitud / Longitud / Tipo) - no change needed.

# --- Data row (row 2) ---
# A2: "Pablo Pinto" -> "Pedro"
$ws.Range("A2").Value = "Pedro"

# B2 holds a mailto: hyperlink; only the displayed text changes here -
# the underlying hyperlink target/relationship is left untouched.
$ws.Range("B2").Value = "pedro@gmail.com"

# C2: "59687412O" -> "56897412M"
$ws.Range("C2").Value = "56897412M"

# D2: 156.26 -> 45.268000000000001
$ws.Range("D2").Value = 45.268000000000001

# E2: -10.265000000000001 -> 56.26
$ws.Range("E2").Value = 56.26

# Selection moves from E2 to A3
[void]$ws.Range("A3").Select()
